$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "26.914.12"
Set-TextValue $ws.Range("E2") "  +0.16%  "
Set-TextValue $ws.Range("D3") "1.545.77"
Set-TextValue $ws.Range("E3") "  -1.21%  "
Set-TextValue $ws.Range("D5") "205.43"
Set-TextValue $ws.Range("E5") "  -0.26%  "
Set-TextValue $ws.Range("E6") "  -0.16%  "
Set-TextValue $ws.Range("E7") "  +0.30%  "
Set-TextValue $ws.Range("E8") "  -0.02%  "
Set-TextValue $ws.Range("D9") "21.32"
Set-TextValue $ws.Range("E9") "  -2.30%  "
Set-TextValue $ws.Range("E10") "  -0.56%  "
Set-TextValue $ws.Range("E11") "  -0.88%  "
Set-TextValue $ws.Range("D12") "1.765.99"
Set-TextValue $ws.Range("E12") "  -1.13%  "
Set-TextValue $ws.Range("D13") "1.542.82"
Set-TextValue $ws.Range("E13") "  -1.34%  "
Set-TextValue $ws.Range("E14") "  -0.98%  "
Set-TextValue $ws.Range("E15") "  -0.65%  "
Set-TextValue $ws.Range("D16") "26.882.15"
Set-TextValue $ws.Range("E16") "  +0.03%  "
Set-TextValue $ws.Range("E17") "  +0.33%  "
Set-TextValue $ws.Range("D18") "213.49"
Set-TextValue $ws.Range("E18") "  -0.83%  "
Set-TextValue $ws.Range("E19") "  -0.02%  "
Set-TextValue $ws.Range("E20") "  -2.56%  "
Set-TextValue $ws.Range("E21") "  +0.27%  "
Set-TextValue $ws.Range("E22") "  -2.84%  "
Set-TextValue $ws.Range("E23") "  -0.33%  "
Set-TextValue $ws.Range("E24") "  -3.58%  "
Set-TextValue $ws.Range("D25") "153.13"
Set-TextValue $ws.Range("E25") "  -0.18%  "
Set-TextValue $ws.Range("E26") "  -1.26%  "
Set-TextValue $ws.Range("E27") "  -0.92%  "
Set-TextValue $ws.Range("E28") "  +0.31%  "
Set-TextValue $ws.Range("E29") "  -0.16%  "
Set-TextValue $ws.Range("E30") "  -1.87%  "
Set-TextValue $ws.Range("E31") "  -1.14%  "
Set-TextValue $ws.Range("E32") "  +1.63%  "
Set-TextValue $ws.Range("D33") "1.360.84"
Set-TextValue $ws.Range("E33") "  -2.91%  "
Set-TextValue $ws.Range("E34") "  +0.33%  "
Set-TextValue $ws.Range("E35") "  +0.08%  "
Set-TextValue $ws.Range("D36") "0.971"
Set-TextValue $ws.Range("E36") "  +6.68%  "
Set-TextValue $ws.Range("E37") "  +0.30%  "
Set-TextValue $ws.Range("E38") "  +0.34%  "
Set-TextValue $ws.Range("D39") "0.517"
Set-TextValue $ws.Range("E39") "  -1.86%  "
Set-TextValue $ws.Range("D40") "0.804"
Set-TextValue $ws.Range("E40") "  -1.12%  "
Set-TextValue $ws.Range("E41") "  +0.30%  "
Set-TextValue $ws.Range("E42") "  -0.18%  "
Set-TextValue $ws.Range("E43") "  -0.71%  "
Set-TextValue $ws.Range("E44") "  +1.49%  "
Set-TextValue $ws.Range("D45") "63.32"
Set-TextValue $ws.Range("E45") "  -0.49%  "
Set-TextValue $ws.Range("E46") "  -3.19%  "
Set-TextValue $ws.Range("D47") "1.680.81"
Set-TextValue $ws.Range("E47") "  -1.16%  "
Set-TextValue $ws.Range("E48") "  -0.57%  "
Set-TextValue $ws.Range("E49") "  +0.66%  "
Set-TextValue $ws.Range("D50") "0.0₇0964"
Set-TextValue $ws.Range("E50") "  -0.93%  "
Set-TextValue $ws.Range("D51") "0.0946"
Set-TextValue $ws.Range("E51") "  -0.41%  "
